# Insert a new weekly price record for "Vega Monumental Concepción - Espinaca"
# at row 102, pushing the existing rows 102:131 down to 103:132.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 102 (shifts 102:131 -> 103:132, dimension grows to R132).
$ws.Rows("102:102").Insert()

# Populate the new row with the latest weekly record (fixed columns repeat the
# values shared by every row in this Espinaca / Vega Monumental Concepción block).
$ws.Cells.Item(102, 1).Value = 11
$ws.Cells.Item(102, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(102, 3).Value = "Bíobío"
$ws.Cells.Item(102, 4).Value = 45135
$ws.Cells.Item(102, 5).Value = 8
$ws.Cells.Item(102, 6).Value = 100112012
$ws.Cells.Item(102, 7).Value = "Espinaca"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 40
$ws.Cells.Item(102, 11).Value = 6500
$ws.Cells.Item(102, 12).Value = 7000
$ws.Cells.Item(102, 13).Value = 6750
$ws.Cells.Item(102, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(102, 15).Value = "Región Metropolitana"
$ws.Cells.Item(102, 16).Value = 675
$ws.Cells.Item(102, 17).Value = 10
$ws.Cells.Item(102, 18).Value = "Hortaliza"
